$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 80.69829787234042
$ws.Range("C2").Value = 88.80627586206897
$ws.Range("D2").Value = 5.71111875
$ws.Range("F2").Value = 82.99227586206896
$ws.Range("H2").Value = 519939466818.5655
$ws.Range("I2").Value = 14342.826620689655
$ws.Range("L2").Value = 21850737.925517242
$ws.Range("M2").Value = 60.597793103448275
$ws.Range("N2").Value = 831735.3982068965
$ws.Range("U2").Value = 211.48862068965516
$ws.Range("X2").Value = 7.415358620689656
$ws.Range("Y2").Value = 87.40430555555555
$ws.Range("B3").Value = 16.198157180766724
$ws.Range("C3").Value = 27.14658726108864
$ws.Range("D3").Value = 4.195846640666898
$ws.Range("F3").Value = 27.38964127737365
$ws.Range("H3").Value = 1962079707344.4812
$ws.Range("I3").Value = 20030.356747389586
$ws.Range("L3").Value = 77389439.40222323
$ws.Range("M3").Value = 10.085351785974003
$ws.Range("N3").Value = 2076664.0604624643
$ws.Range("U3").Value = 690.6488144212989
$ws.Range("X3").Value = 5.402807340091799
$ws.Range("Y3").Value = 16.245341513811987
$ws.Range("B5").Value = 80.69829787234042
$ws.Range("B6").Value = 80.69829787234042
$ws.Range("M6").Value = 60.675
$ws.Range("B9").Value = 12.780011999999996
$ws.Range("M9").Value = 7.798475999999996
$ws.Range("B10").Value = 13.436702127659572
$ws.Range("B11").Value = 0.20072489269093607
$ws.Range("C11").Value = 0.30568320760631645
$ws.Range("D11").Value = 0.7346803357340272
$ws.Range("F11").Value = 0.3300263909245546
$ws.Range("H11").Value = 3.7736694991633613
$ws.Range("I11").Value = 1.396541788945257
$ws.Range("L11").Value = 3.5417311610262834
$ws.Range("M11").Value = 0.16643100795365603
$ws.Range("N11").Value = 2.4967845121651155
$ws.Range("U11").Value = 3.2656547296451377
$ws.Range("X11").Value = 0.7285969049450124
$ws.Range("Y11").Value = 0.18586431652942076
$ws.Range("B12").Value = -1.2918474361853092
$ws.Range("C12").Value = -2.5225506749590334
$ws.Range("D12").Value = 0.2994784892801578
$ws.Range("F12").Value = -1.450059970537098
$ws.Range("H12").Value = 7.4781596907236745
$ws.Range("I12").Value = 2.0375608312488693
$ws.Range("L12").Value = 7.8660128983431346
$ws.Range("M12").Value = -0.055871291631802225
$ws.Range("N12").Value = 4.774058858157279
$ws.Range("U12").Value = 9.548374445556385
$ws.Range("X12").Value = 1.3320562133951013
$ws.Range("Y12").Value = -1.3217050379494302
$ws.Range("B14").Value = 1.5447454302717345
$ws.Range("C14").Value = 4.929489407085548
$ws.Range("D14").Value = -0.9770089904580632
$ws.Range("F14").Value = 0.6630652679147837
$ws.Range("H14").Value = 62.23216894296577
$ws.Range("I14").Value = 4.174797097366903
$ws.Range("L14").Value = 68.00581013047517
$ws.Range("M14").Value = 0.13802317780256645
$ws.Range("N14").Value = 25.881556055086705
$ws.Range("U14").Value = 101.17800774626551
$ws.Range("X14").Value = 1.2250453778464818
$ws.Range("Y14").Value = 0.5218438779401575
